$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: WELCOME
$ws.Range("B2").Value = 'Willkommen zum Dreiklangsfolgen-Test'
$ws.Range("C2").Value = 'Welcome to the Harmony Progression Discrimination Test'

# Row 3: INTRO_TEXT
$ws.Range("B3").Value = 'Mit diesem Test wollen wir herausfinden, wie unsere Ohren Dreiklänge wahrnehmen. Ein Dreiklang ist eine Gruppe von drei Tönen, die gleichzeitig abgespielt werden. \\ Bei jeder Frage wirst du zwei kurze Dreiklangsfolgen hören, die aus jeweils vier Dreiklängen bestehen. Die beiden Folgen sind exakt dieselben bis auf einen einzelnen Dreiklang. Deine Aufgabe ist es herauszufinden, welcher Dreiklang anders ist, indem du auf die entsprechende Zahl klickst. Die Zahlen passen zu der Reihenfolge der Dreiklänge. \\  Bitte höre dir beide Dreiklangsfolgen immer ganz an, bevor du deine Antwort auswählst. Falls du dir nicht sicher bist, wähle einfach deine beste Vermutung aus.  '
$ws.Range("C3").Value = 'With this test we want to investigate how we understand chords. A chord is a group of musical notes played at the same time. \\  For each question on this test, you will hear two short chord progressions consisting of four chords each. The two progressions will be identical with the exception of one chord. Your task is to determine which chord is different by selecting the appropriate number between 1 and 4. These numbers correspond to the order of chords presented.  \\ Please listen to each sound clip in full before making your decision. If you don’t know the answer, give your best guess.'

# Row 6: ITEM_INSTRUCTION
$ws.Range("B6").Value = 'Bitte höre dir die folgenden Musikausschnitte an und entscheide, welcher Dreiklang verändert wurde. Wähle die entsprechende Zahl zwischen 1 und 4. Die Zahlen stimmen mit der Reihenfolge der gehörten Dreiklänge überein.'
$ws.Range("C6").Value = 'Please listen to the following clips and select which chord was different. Select the appropriate number between 1 and 4. These numbers correspond to the order of chords presented.'

# Row 9: FEEDBACK
$ws.Range("B9").Value = 'Du hast **{{accuracy}} %** der veränderten Dreiklänge richtig erkannt.'
$ws.Range("C9").Value = 'You recognized **{{accuracy}} %**  of the altered chords correctly.'

# Row 13: SUCCESS
$ws.Range("B13").Value = 'Du hast den Dreiklangsfolgen-Test erfolgreich beendet.'
$ws.Range("C13").Value = 'You have completed the Harmony Progression Discrimination Test.'

# Row 14: TESTNAME
$ws.Range("B14").Value = 'Dreiklangsfolgen-Test'
$ws.Range("C14").Value = 'Harmony Progression Discrimination Test'

# Row 19: INSTRUCTIONS
$ws.Range("B19").Value = 'Zuerst wirst du Beispiele hören und dann ein paar Übungsaufgaben machen.'
$ws.Range("C19").Value = 'Try a couple of practice questions before the test begins.'

# Row 20: SAMPLE1a
$ws.Range("B20").Value = 'Zuerst wirst du Beispiele hören und dann ein paar Übungsaufgaben machen. \\ Bitte höre dir die folgenden Musikausschnitte an und entscheide, welcher Dreikläng verändert wurde. Überlege dir welche Zahle zwischen 1 und 4 zu dem veränderten Dreiklänge gehört. Die Zahlen stimmen mit der Reihenfolge der gehörten Dreiklänge überein.'
$ws.Range("C20").Value = 'Try a couple of practice questions before the test begins. \\ Please listen to the following clips and select which chord was different. Think about which number between 1 and 4 fits the different chord. These numbers correspond to the order of chords presented.'

# Row 21: SAMPLE1b
$ws.Range("B21").Value = 'In diesem Beispiel wurde der dritte Dreiklang verändert. Die richtige Antwort wäre also **Nummer 3**. Es folgen nun zwei Übungsfragen.'
$ws.Range("C21").Value = ' Here, the third chord was different, so the correct answer would have been **number 3**. Now you will see two practice questions.'

# Row 22: PRACTICE1
$ws.Range("B22").Value = '**Übungsfrage 1** \\ Bitte höre dir die folgenden Musikausschnitte an und entscheide, welcher Dreiklang verändert wurde. Überlege dir welche Zahle zwischen 1 und 4 zu dem veränderten Dreiklänge gehört. Die Zahlen stimmen mit der Reihenfolge der gehörten Dreiklänge überein.'
$ws.Range("C22").Value = '**Practice question 1** \\ Please listen to the following clips and select which chord was different. Think about which number between 1 and 4 fits the differenct chord. These numbers correspond to the order of chords presented.'

# Row 23: PRACTICE2
$ws.Range("B23").Value = '**{{feedback}}** \\ Hier ist ein weiteres Beispiel. \\ **Übungsfrage 2:** \\ Bitte höre dir die folgenden Musikausschnitte an und entscheide, welcher Dreiklang verändert wurde. Überlege dir welche Zahle zwischen 1 und 4 zu dem veränderten Dreiklang gehört. Die Zahlen stimmen mit der Reihenfolge der gehörten Dreiklänge überein.'
$ws.Range("C23").Value = ' **{{feedback}}** \\ Here’s another example. \\ **Practice question 2:** \\ Please listen to the following clips and select which chord was different. Think about which number between 1 and 4 fits the differenct chord. These numbers correspond to the order of chords presented.'

# Row 28: RESULTS_SAVED
$ws.Range("B28").Value = 'Deine Ergebnisse wurden gespeichert.'
$ws.Range("C28").Value = 'Your results have been saved.'

# Row 29: CLOSE_BROWSER
$ws.Range("B29").Value = 'Du kannst den Browsertab jetzt schließen.'
$ws.Range("C29").Value = 'You may now close the browser tab.'

# Update view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C22").Select()
